# Created method for partners dropdown in loanentry and loanreview screen
#
# 1) Rename the existing sheet "LOGINCREDS" -> "logincrds"
# 2) Insert a new "partners" sheet right after it, listing the partner
#    names that back the dropdown, with a "name" header cell.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "logincrds"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "partners"

$partners = @(
    "sonata",
    "Samasta",
    "Pahal",
    "Cashpor",
    "Magalir",
    "loadtestingpartner11",
    "loadtestingpartner12",
    "loadtestingpartner13",
    "loadtestingpartner14",
    "loadtestingpartner15",
    "Sugmya",
    "1M2S3M",
    "subhlakshmi",
    "CTL",
    "Samavesh PR",
    "Svasti Microfinance",
    "Mitrata",
    "SWARA FINCARE LIMITED",
    "Seeds Fincap",
    "Subhlakshmi Finance Pvt Ltd",
    "Maximal",
    "Kiara",
    "IREP",
    "Midland_test",
    "Sonata PR",
    "Magalir PR",
    "sona1ta1",
    "Pahal PR",
    "Sugmya PR",
    "MSM",
    "subhlakshmi PR",
    "Capital Trust Limited PR",
    "Svasti PR",
    "samavesh",
    "Mitrata PR",
    "MIDLAND_MICROFIN_LIMITED",
    "MAXIMAL-KISCORE"
)

# Write all of the partner names first (rows 2..38) so the shared-string
# table picks them up in this order, then write the "name" header last.
for ($i = 0; $i -lt $partners.Length; $i++) {
    $ws2.Cells.Item($i + 2, 1).Value = $partners[$i]
}
$ws2.Cells.Item(1, 1).Value = "name"

# Build the Courier New / vertical-centered look once on the first data
# cell, then fan it out to the rest of the list via a format-only paste
# so we don't re-derive the style on every row.
$first = $ws2.Range("A2")
$first.Font.Name = "Courier New"
$first.Font.Family = 3
$first.Font.Size = 10
$first.Font.Color = 0
$first.VerticalAlignment = -4108

$first.Copy()
$ws2.Range("A3:A38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header gets a green highlight fill.
$header = $ws2.Range("A1")
$header.Interior.Color = 5287936

$ws2.Columns.Item(1).ColumnWidth = 31.08984375

$ws2.Range("A1").Select()

$ws1.Activate()
$ws1.Range("M20").Select()
